$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F16").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F18").Value = "application instructions || env warning - species"
$ws.Range("F19").Value = "pollinator"
$ws.Range("F20").Value = "32_physical_and_chemical_hazards"
$ws.Range("F23").Value = "135_product_information"
$ws.Range("F24").Value = "application instructions"
$ws.Range("F25").Value = "mixing"
$ws.Range("F26").Value = "mixing"
$ws.Range("F27").Value = "mixing"
$ws.Range("F28").Value = "application instructions"
$ws.Range("F29").Value = "application instructions"
$ws.Range("F30").Value = "application instructions"
$ws.Range("F31").Value = "application instructions"
$ws.Range("F32").Value = "use restrictions"
$ws.Range("F33").Value = "use restrictions"
$ws.Range("F34").Value = "application instructions"
$ws.Range("F44").Value = "application instructions"
$ws.Range("F45").Value = "application instructions"
$ws.Range("F46").Value = "application instructions"
$ws.Range("F47").Value = "application instructions"
$ws.Range("F48").Value = "application instructions"
$ws.Range("F49").Value = "application instructions"
$ws.Range("F50").Value = "application instructions"
$ws.Range("F51").Value = "application instructions"
$ws.Range("F52").Value = "mixing || application instructions"
$ws.Range("F53").Value = "application instructions"
$ws.Range("F54").Value = "application instructions"
$ws.Range("F55").Value = "application instructions"
$ws.Range("F82").Value = "application instructions"
$ws.Range("F84").Value = "application instructions"
$ws.Range("F85").Value = "application instructions"
$ws.Range("F86").Value = "application instructions"
$ws.Range("F87").Value = "application instructions"
$ws.Range("F88").Value = "application instructions"
$ws.Range("F92").Value = "use restrictions || irrigation"
$ws.Range("F93").Value = "154_pesticide_storage"
